$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:J2").Value = "unknown"
